# إضافة حدث جديد في Card20 by admin at 2025-12-08 09:03:26
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# Row 15 currently has empty inlineStr cells in B:K and M; fill them with "nan"
# to match the exported pandas-style placeholder used throughout this sheet.
$ws.Range("B15:K15").Value = "nan"
$ws.Range("M15").Value = "nan"

# Add the new service event as row 16
$ws.Range("A16").Value = "20"
$ws.Range("L16").Value = "30\9\2024"
$ws.Range("N16").Value = "تم تاكيد المعيار"
$ws.Range("O16").Value = "الخبير"
